$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.439.15"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").Value = "1.842.82"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  +1.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.19"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.013"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4761"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07464"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.17%  "
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "1.850.29"
$ws.Range("E12").Value = "  +1.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07391"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.485"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.49"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.601"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008850"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.015"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").Value = "27.456.34"
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").Value = "2.080.43"
$ws.Range("E24").Value = "  +2.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.904"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.56"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.66"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.166"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.294"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.11"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08986"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7606"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.568"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.953"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.015"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.107"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05363"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01970"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.006"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.322"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5361"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.375"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.93%  "
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.558"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4980"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.65"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.39%  "
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.34"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.687"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06326"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.30%  "
